# Chapter 3 is added
#
# Inserts three new paragraphs (an empty spacer paragraph, the "## Глава 3:
# Вход в систему" heading paragraph, and its body paragraph) right after
# the "Основные понятия и команды Git." paragraph that closes Chapter 2,
# and right before the trailing empty paragraph that precedes the
# section break.

$d = $word.ActiveDocument

# Locate the paragraph that ends Chapter 2 ("Основные понятия и команды
# Git.") by searching for a distinctive substring, so the script does not
# depend on hard-coded paragraph indices.
$searchRange = $d.Content
$found = $searchRange.Find.Execute("Основные понятия и команды", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) {
    throw "edit.ps1: could not find the Chapter 2 anchor paragraph"
}

$anchorIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $candidate = $d.Paragraphs.Item($i)
    if ($candidate.Range.Start -lt $searchRange.End -and $candidate.Range.End -ge $searchRange.End) {
        $anchorIndex = $i
    }
}
if ($anchorIndex -eq -1) {
    throw "edit.ps1: could not resolve the Chapter 2 anchor paragraph index"
}

$anchorParagraph = $d.Paragraphs.Item($anchorIndex)
$insertParagraph = $anchorParagraph.Next()
if ($insertParagraph) {
    $insertRange = $insertParagraph.Range
    $insertRange.Collapse(1)   # wdCollapseStart - insert before the trailing empty paragraph
} else {
    # Anchor was the last paragraph in the story; insert right after it instead.
    $insertRange = $anchorParagraph.Range
    $insertRange.Collapse(0)   # wdCollapseEnd
}

$xmlFragment = @"
<?xml version="1.0" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
<w:p><w:pPr><w:spacing w:after="0"/><w:ind w:firstLine="709"/><w:jc w:val="both"/></w:pPr></w:p>
<w:p><w:pPr><w:spacing w:after="0"/><w:ind w:firstLine="709"/><w:jc w:val="both"/></w:pPr><w:r><w:t>## Глава 3: Вход в систему</w:t></w:r></w:p>
<w:p><w:pPr><w:spacing w:after="0"/><w:ind w:firstLine="709"/><w:jc w:val="both"/></w:pPr><w:r><w:t>Раздел по новой функциональности входа в систему.</w:t></w:r></w:p>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
"@

$insertRange.InsertXML($xmlFragment)

Write-Output "Chapter 3 inserted; paragraph count is now $($d.Paragraphs.Count)"
